$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): new columns N ("ТК_оригинал") and P ("на момент выгрузки в элжуре") ---
$ws.Range("N2").Value = "ТК_оригинал"
$ws.Range("P2").Value = "на момент выгрузки в элжуре"

# --- Remove the stray "Хочу 4" markers that used to live in column O for rows 4, 7, 11 ---
$ws.Range("O4").ClearContents()
$ws.Range("O7").ClearContents()
$ws.Range("O11").ClearContents()
# (O28 keeps its existing "Хочу 4" value - left untouched)

# --- New column P: copy of column N ("ТК_оригинал" control values) for every student row ---
for ($r = 4; $r -le 31; $r++) {
    $nVal = $ws.Range("N$r").Value2
    $ws.Range("P$r").Value = $nVal
}

# Format column P (rows 4-31) with the same thick box border / centred wrapped text
# used elsewhere in the sheet - apply to the whole block in one shot so the engine
# settles on a single consolidated cell style.
$pRange = $ws.Range("P4:P31")
$pRange.Borders.LineStyle = 1
$pRange.Borders.Weight = 4
$pRange.Borders.Color = 0
$pRange.HorizontalAlignment = -4108
$pRange.VerticalAlignment = -4108
$pRange.WrapText = $true

# --- New column Q: check formula P-N (should evaluate to 0 everywhere) ---
$ws.Range("Q4").Formula = "=P4-N4"
$ws.Range("Q5:Q31").Formula = "=P5-N5"
